# Updated cryptos list (price + 1h volume change columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.798.65"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.349.44"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").Value = "3.343.07"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.173"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.64%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.637"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.52"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "3.883.65"
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "3.339.34"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("D19").Value = "64.648.74"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.988"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("E24").Value = "  -0.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.73%  "
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "62.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "574.37"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "0.0₃0743"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.370"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("D42").Value = "3.087.48"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("E45").Value = "  +3.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.58"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("E50").Value = "  -2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.44%  "
